$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @{ B=10949.0256794104;  C=10138.9613257071;  E=6844.01555828536;  F=-3.38679650031266 }
  3  = @{ B=11027.8441277067;  C=10172.7327820601;  E=6908.60306675794;  F=243.711493700752 }
  4  = @{ B=10877.3098276712;  C=10071.2453393396;  E=6801.39799789617;  F=235.015972384822 }
  5  = @{ B=10681.5020285746;  C=9436.64328802674;  E=6656.73977179742;  F=202.54679415934 }
  6  = @{ B=8330.02363792758;  C=8854.2971956287;   E=6964.39251149521;  F=191.10123779683 }
  7  = @{ B=8231.05518265296;  C=8710.56234215091;  E=6952.36995051702;  F=184.61134552783 }
  8  = @{ B=8386.79580737583;  C=8729.99027894652;  E=7376.64878490783;  F=203.099127660598 }
  9  = @{ B=9707.03874689317;  C=9731.17034323013;  E=7791.81093033999;  F=262.113386398755 }
  10 = @{ B=9707.03874689317;  C=9407.43104650919;  E=7791.81093033999;  F=248.624249035382 }
  11 = @{ B=9707.03874689317;  C=9402.13697322575;  E=7791.81093033999;  F=248.403662648572 }
  12 = @{ B=9707.03874689317;  C=8660.66086194875;  E=7791.81093033999;  F=217.508824678697 }
  13 = @{ B=8488.36206656772;  C=8448.79064244595;  E=7389.01917685293;  F=191.897909137453 }
  14 = @{ B=8386.79580737583;  C=8392.87341934355;  E=7376.72537237786;  F=189.055782988392 }
  15 = @{ B=9746.09605670163;  C=9115.19056643991;  E=8065.19356806475;  F=247.838505604361 }
}

foreach ($row in $data.Keys) {
  $vals = $data[$row]
  $ws.Range("B$row").Value = $vals.B
  $ws.Range("C$row").Value = $vals.C
  $ws.Range("E$row").Value = $vals.E
  $ws.Range("F$row").Value = $vals.F
}
